$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(1).RowHeight = 25.5
Write-Output "done"
